$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2944.15
$ws.Range("I111").Value = 3233.625
$ws.Range("K111").Value = 9700.875
$ws.Range("M111").Value = -6633.875

$ws.Range("H132").Value = 3358.025
$ws.Range("I132").Value = 3477.162
$ws.Range("K132").Value = 10431.486
$ws.Range("M132").Value = -7901.485999999999

$ws.Range("H138").Value = 2994.8171
$ws.Range("I138").Value = 1910.591
$ws.Range("J138").Value = 3392.3667
$ws.Range("K138").Value = 5731.772999999999
$ws.Range("L138").Value = 10177.1001
$ws.Range("M138").Value = -591.7729999999992
$ws.Range("N138").Value = -20457.1001


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5710.52
$ws.Range("I32").Value = 3311.2954
$ws.Range("J32").Value = 23304.834
$ws.Range("K32").Value = 3311.2954
$ws.Range("L32").Value = 23304.834
$ws.Range("M32").Value = -3024.2954
$ws.Range("N32").Value = -23878.834

$ws.Range("H61").Value = 7859
$ws.Range("I61").Value = 9870.223
$ws.Range("K61").Value = 9870.223
$ws.Range("M61").Value = -9658.223

$ws.Range("H74").Value = 27818.414
$ws.Range("I74").Value = 5911.143
$ws.Range("K74").Value = 5911.143
$ws.Range("M74").Value = -5037.143

$ws.Range("H77").Value = 27818.414
$ws.Range("I77").Value = 5911.143
$ws.Range("K77").Value = 29555.715
$ws.Range("M77").Value = -25187.715

$ws.Range("H122").Value = 3217.9395
$ws.Range("I122").Value = 1743.0476
$ws.Range("J122").Value = 5799
$ws.Range("K122").Value = 5229.142800000001
$ws.Range("L122").Value = 17397
$ws.Range("M122").Value = -2779.142800000001
$ws.Range("N122").Value = -22297

$ws.Range("H127").Value = 99755.55499999999
$ws.Range("J127").Value = 99755.55499999999
$ws.Range("L127").Value = 99755.55499999999
$ws.Range("N127").Value = -109675.555

$ws.Range("H132").Value = 36073.965
$ws.Range("I132").Value = 9464.154
$ws.Range("J132").Value = 56422.65
$ws.Range("K132").Value = 28392.462
$ws.Range("L132").Value = 169267.95
$ws.Range("M132").Value = -25862.462
$ws.Range("N132").Value = -174327.95

$ws.Range("H136").Value = 7859
$ws.Range("I136").Value = 9870.223
$ws.Range("K136").Value = 29610.669
$ws.Range("M136").Value = -27060.669


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 422
$ws.Range("I22").Value = 422
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 422
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -249
$ws.Range("N22").ClearContents()

$ws.Range("H111").Value = 60700
$ws.Range("J111").Value = 60700
$ws.Range("L111").Value = 60700
$ws.Range("N111").Value = -68880


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1639.9231
$ws.Range("I16").Value = 1393.6364
$ws.Range("J16").Value = 2994.5
$ws.Range("K16").Value = 1393.6364
$ws.Range("L16").Value = 2994.5
$ws.Range("M16").Value = -1106.6364
$ws.Range("N16").Value = -3568.5

$ws.Range("H31").Value = 30732.154
$ws.Range("I31").Value = 37166.668
$ws.Range("K31").Value = 37166.668
$ws.Range("M31").Value = -36871.668

$ws.Range("H34").Value = 30732.154
$ws.Range("I34").Value = 37166.668
$ws.Range("K34").Value = 37166.668
$ws.Range("M34").Value = -36964.668

$ws.Range("H58").Value = 11272
$ws.Range("I58").Value = 16051
$ws.Range("J58").Value = 4581.4
$ws.Range("K58").Value = 16051
$ws.Range("L58").Value = 4581.4
$ws.Range("M58").Value = -15848
$ws.Range("N58").Value = -4987.4

$ws.Range("H113").Value = 1639.9231
$ws.Range("I113").Value = 1393.6364
$ws.Range("J113").Value = 2994.5
$ws.Range("K113").Value = 1393.6364
$ws.Range("L113").Value = 2994.5
$ws.Range("M113").Value = 776.3635999999999
$ws.Range("N113").Value = -7334.5

$ws.Range("H136").Value = 11272
$ws.Range("I136").Value = 16051
$ws.Range("J136").Value = 4581.4
$ws.Range("K136").Value = 48153
$ws.Range("L136").Value = 13744.2
$ws.Range("M136").Value = -45603
$ws.Range("N136").Value = -18844.2

$ws.Range("H141").Value = 123195.64
$ws.Range("J141").Value = 131647.83
$ws.Range("L141").Value = 131647.83
$ws.Range("N141").Value = -142007.83


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7372
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 18000
$ws.Range("N65").Value = -24864
$ws.Range("M65").ClearContents()

$ws.Range("H74").Value = 11224.111
$ws.Range("I74").Value = 5339
$ws.Range("J74").Value = 14166.667
$ws.Range("K74").Value = 16017
$ws.Range("L74").Value = 42500.001
$ws.Range("M74").Value = -14956
$ws.Range("N74").Value = -44622.001

$ws.Range("H77").Value = 11224.111
$ws.Range("I77").Value = 5339
$ws.Range("J77").Value = 14166.667
$ws.Range("K77").Value = 48051
$ws.Range("L77").Value = 127500.003
$ws.Range("M77").Value = -42747
$ws.Range("N77").Value = -138108.003

$ws.Range("H81").Value = 6873.1665
$ws.Range("J81").Value = 6873.1665
$ws.Range("L81").Value = 20619.4995
$ws.Range("N81").Value = -22865.4995

$ws.Range("H84").Value = 6873.1665
$ws.Range("J84").Value = 6873.1665
$ws.Range("L84").Value = 61858.4985
$ws.Range("N84").Value = -73090.4985


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 289636
$ws.Range("I122").Value = 357788.75
$ws.Range("J122").Value = 5666.1665
$ws.Range("K122").Value = 1073366.25
$ws.Range("L122").Value = 16998.4995
$ws.Range("M122").Value = -1070916.25
$ws.Range("N122").Value = -21898.4995

$ws.Range("H132").Value = 9058.333000000001
$ws.Range("J132").Value = 11437.154
$ws.Range("L132").Value = 34311.462
$ws.Range("N132").Value = -39371.462


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5224.385
$ws.Range("I7").Value = 4384.1665
$ws.Range("J7").Value = 5944.5713
$ws.Range("K7").Value = 4384.1665
$ws.Range("L7").Value = 5944.5713
$ws.Range("M7").Value = -4272.1665
$ws.Range("N7").Value = -6168.5713

$ws.Range("H40").Value = 4674.0454
$ws.Range("I40").Value = 3043.7368
$ws.Range("K40").Value = 3043.7368
$ws.Range("M40").Value = -2907.7368

$ws.Range("H61").Value = 2400
$ws.Range("I61").Value = 2400
$ws.Range("K61").Value = 2400
$ws.Range("M61").Value = -2198

$ws.Range("H82").Value = 1939.8214
$ws.Range("I82").Value = 2415.125
$ws.Range("K82").Value = 2415.125
$ws.Range("M82").Value = -2054.125

$ws.Range("H85").Value = 1939.8214
$ws.Range("I85").Value = 2415.125
$ws.Range("K85").Value = 2415.125
$ws.Range("M85").Value = -1167.125

$ws.Range("H93").Value = 6706.1816
$ws.Range("I93").Value = 2124.6
$ws.Range("J93").Value = 16523.857
$ws.Range("K93").Value = 2124.6
$ws.Range("L93").Value = 16523.857
$ws.Range("M93").Value = -876.5999999999999
$ws.Range("N93").Value = -19019.857

$ws.Range("H113").Value = 2400
$ws.Range("I113").Value = 2400
$ws.Range("K113").Value = 2400
$ws.Range("M113").Value = -230

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H126").Value = 5224.385
$ws.Range("I126").Value = 4384.1665
$ws.Range("J126").Value = 5944.5713
$ws.Range("K126").Value = 13152.4995
$ws.Range("L126").Value = 17833.7139
$ws.Range("M126").Value = -10682.4995
$ws.Range("N126").Value = -22773.7139

$ws.Range("H136").Value = 96864.87
$ws.Range("I136").Value = 171002.5
$ws.Range("K136").Value = 513007.5
$ws.Range("M136").Value = -510457.5


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8307.75
$ws.Range("J74").Value = 8307.75
$ws.Range("L74").Value = 8307.75
$ws.Range("N74").Value = -10179.75

$ws.Range("H77").Value = 8307.75
$ws.Range("J77").Value = 8307.75
$ws.Range("L77").Value = 24923.25
$ws.Range("N77").Value = -34283.25

$ws.Range("H96").Value = 3818.7334
$ws.Range("I96").Value = 2767
$ws.Range("J96").Value = 4344.6
$ws.Range("K96").Value = 2767
$ws.Range("L96").Value = 4344.6
$ws.Range("M96").Value = -1394
$ws.Range("N96").Value = -7090.6

$ws.Range("H100").Value = 2373.652
$ws.Range("I100").Value = 2742.8948
$ws.Range("K100").Value = 5485.7896
$ws.Range("M100").Value = -4944.7896

$ws.Range("H122").Value = 3243.125
$ws.Range("I122").Value = 1922.3914
$ws.Range("K122").Value = 5767.174199999999
$ws.Range("M122").Value = -3317.174199999999


Write-Host "All updates applied."